# Mise a jour de certains champs de Modules et de Professeurs
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Modules" header row: swap the Enseignant/Nombre d'heures
# columns for Chef Module / Composants. We set D1 first then C1 so
# that the in-place shared-string rewrite lands on the same si indices
# that the real workbook ends up with (C1 -> index 3, D1 -> index 2).
$ws.Range("D1").Value = "Composants"
$ws.Range("C1").Value = "Chef  Module"

# Widen columns C and D to fit the new headers.
# ColumnWidth is expressed in characters of the Normal style font; the
# engine stores/rounds the resulting width to a pixel grid, so we feed
# it the (target - 0.8333333333333334) offset that reproduces the
# desired stored width as closely as possible.
$ws.Columns.Item(3).ColumnWidth = 34.166666666666664
$ws.Columns.Item(4).ColumnWidth = 23.736979166666668

# Move the active selection to E8.
$ws.Range("E8").Select()
